$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force the "Price" column cells to Text format before assigning, so that
# numeric-looking strings (e.g. "544.51") are stored as text, matching the
# original inline-string cell content instead of being parsed as numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '60.935.04'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '2.370.15'
$ws.Range('E3').Value = '  -3.75%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '544.51'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('D6').Value = '140.46'
$ws.Range('E6').Value = '  -2.90%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.533'
$ws.Range('E8').Value = '  -11.64%  '
$ws.Range('D9').Value = '2.372.07'
$ws.Range('E9').Value = '  -3.62%  '
$ws.Range('E10').Value = '  -2.01%  '
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('D12').Value = '5.24'
$ws.Range('E12').Value = '  -3.14%  '
$ws.Range('D13').Value = '0.342'
$ws.Range('E13').Value = '  -3.29%  '
$ws.Range('D14').Value = '25.36'
$ws.Range('E14').Value = '  -1.87%  '
$ws.Range('D15').Value = '2.801.10'
$ws.Range('E15').Value = '  -3.60%  '
$ws.Range('D16').Value = '0.0000166'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').Value = '60.669.64'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').Value = '2.374.45'
$ws.Range('E18').Value = '  -3.39%  '
$ws.Range('D19').Value = '10.72'
$ws.Range('E19').Value = '  -3.06%  '
$ws.Range('D20').Value = '4.11'
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').Value = '317.10'
$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('E22').Value = '  -3.56%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = '1.91'
$ws.Range('E24').Value = '  +8.30%  '
$ws.Range('D25').Value = '63.22'
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = '0.0₃0940'
$ws.Range('E27').Value = '  -4.39%  '
$ws.Range('D28').Value = '2.492.02'
$ws.Range('E28').Value = '  -3.30%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').Value = '532.30'
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = '7.74'
$ws.Range('E30').Value = '  +2.50%  '
$ws.Range('E31').Value = '  -3.64%  '
$ws.Range('D32').Value = '8.02'
$ws.Range('E32').Value = '  -3.69%  '
$ws.Range('D33').Value = '0.144'
$ws.Range('E33').Value = '  -3.83%  '
$ws.Range('E34').Value = '  -3.36%  '
$ws.Range('E35').Value = '  -0.76%  '
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '5.51'
$ws.Range('E37').Value = '  -4.89%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '4.63'
$ws.Range('E38').Value = '  -3.51%  '
$ws.Range('D39').Value = '0.374'
$ws.Range('E39').Value = '  -1.27%  '
$ws.Range('D40').Value = '17.98'
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').Value = '137.49'
$ws.Range('E43').Value = '  -5.70%  '
$ws.Range('D44').Value = '40.25'
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('D45').Value = '2.24'
$ws.Range('E45').Value = '  -3.71%  '
$ws.Range('D46').Value = '140.08'
$ws.Range('E46').Value = '  -4.49%  '
$ws.Range('D47').Value = '3.56'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('D48').Value = '20.17'
$ws.Range('E48').Value = '  -1.98%  '
$ws.Range('D49').Value = '0.0516'
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('D50').Value = '0.575'
$ws.Range('D51').Value = '0.0908'
$ws.Range('E51').Value = '  -3.26%  '
